$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 15114.875
$ws.Range("I62").Value = 2434.75
$ws.Range("J62").Value = 27795
$ws.Range("K62").Value = 2434.75
$ws.Range("L62").Value = 27795
$ws.Range("M62").Value = -1810.75
$ws.Range("N62").Value = -29043
$ws.Range("H65").Value = 15114.875
$ws.Range("I65").Value = 2434.75
$ws.Range("J65").Value = 27795
$ws.Range("K65").Value = 12173.75
$ws.Range("L65").Value = 138975
$ws.Range("M65").Value = -9053.75
$ws.Range("N65").Value = -145215
$ws.Range("H74").Value = 23816280
$ws.Range("I74").Value = 33338460
$ws.Range("J74").Value = 10833.167
$ws.Range("K74").Value = 33338460
$ws.Range("L74").Value = 10833.167
$ws.Range("M74").Value = -33337524
$ws.Range("N74").Value = -12705.167
$ws.Range("H77").Value = 23816280
$ws.Range("I77").Value = 33338460
$ws.Range("J77").Value = 10833.167
$ws.Range("K77").Value = 166692300
$ws.Range("L77").Value = 54165.835
$ws.Range("M77").Value = -166687620
$ws.Range("N77").Value = -63525.835
$ws.Range("H86").Value = 154324000
$ws.Range("I86").Value = 222222910
$ws.Range("K86").Value = 222222910
$ws.Range("M86").Value = -222221787
$ws.Range("H89").Value = 154324000
$ws.Range("I89").Value = 222222910
$ws.Range("K89").Value = 1111114550
$ws.Range("M89").Value = -1111108934
$ws.Range("H92").Value = 1506.1765
$ws.Range("I92").Value = 677.3077
$ws.Range("J92").Value = 4200
$ws.Range("K92").Value = 677.3077
$ws.Range("L92").Value = 4200
$ws.Range("M92").Value = 570.6923
$ws.Range("N92").Value = -6696
$ws.Range("H137").Value = 2000
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 2447916.2
$ws.Range("J138").Value = 2951292
$ws.Range("L138").Value = 8853876
$ws.Range("N138").Value = -8864156
$ws.Range("H141").Value = 76925060
$ws.Range("I141").Value = 100001224
$ws.Range("J141").Value = 4499.6665
$ws.Range("K141").Value = 300003672
$ws.Range("L141").Value = 13498.9995
$ws.Range("M141").Value = -299998492
$ws.Range("N141").Value = -23858.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4243.909
$ws.Range("I32").Value = 3340.8157
$ws.Range("K32").Value = 3340.8157
$ws.Range("M32").Value = -3053.8157
$ws.Range("H74").Value = 50289.145
$ws.Range("I74").Value = 102626.2
$ws.Range("J74").Value = 2710
$ws.Range("K74").Value = 102626.2
$ws.Range("L74").Value = 2710
$ws.Range("M74").Value = -101752.2
$ws.Range("N74").Value = -4458
$ws.Range("H77").Value = 50289.145
$ws.Range("I77").Value = 102626.2
$ws.Range("J77").Value = 2710
$ws.Range("K77").Value = 513131
$ws.Range("L77").Value = 13550
$ws.Range("M77").Value = -508763
$ws.Range("N77").Value = -22286
$ws.Range("H109").Value = 59386.5
$ws.Range("J109").Value = 59386.5
$ws.Range("L109").Value = 59386.5
$ws.Range("N109").Value = -62160.5
$ws.Range("H111").Value = 59387
$ws.Range("J111").Value = 59387
$ws.Range("L111").Value = 59387
$ws.Range("N111").Value = -67567
$ws.Range("H112").Value = 58886.5
$ws.Range("J112").Value = 58886.5
$ws.Range("L112").Value = 58886.5
$ws.Range("N112").Value = -61840.5
$ws.Range("H114").Value = 64445
$ws.Range("J114").Value = 64445
$ws.Range("L114").Value = 64445
$ws.Range("N114").Value = -73123
$ws.Range("H132").Value = 1618803.6
$ws.Range("I132").Value = 2087747.5
$ws.Range("J132").Value = 10996
$ws.Range("K132").Value = 6263242.5
$ws.Range("L132").Value = 32988
$ws.Range("M132").Value = -6260712.5
$ws.Range("N132").Value = -38048

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 39795
$ws.Range("J2").Value = 39795
$ws.Range("L2").Value = 39795
$ws.Range("N2").Value = -40021
$ws.Range("H94").Value = 29414362
$ws.Range("I94").Value = 40001064
$ws.Range("K94").Value = 40001064
$ws.Range("M94").Value = -40000613
$ws.Range("H105").Value = 2816.3823
$ws.Range("I105").Value = 1821.88
$ws.Range("J105").Value = 5578.8887
$ws.Range("K105").Value = 1821.88
$ws.Range("L105").Value = 5578.8887
$ws.Range("M105").Value = -74.88000000000011
$ws.Range("N105").Value = -9072.8887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 175.90475
$ws.Range("I7").Value = 80.13333
$ws.Range("K7").Value = 80.13333
$ws.Range("M7").Value = 32.86667
$ws.Range("H31").Value = 5856.2573
$ws.Range("I31").Value = 1206.5385
$ws.Range("J31").Value = 8603.817999999999
$ws.Range("K31").Value = 1206.5385
$ws.Range("L31").Value = 8603.817999999999
$ws.Range("M31").Value = -911.5385000000001
$ws.Range("N31").Value = -9193.817999999999
$ws.Range("H34").Value = 5856.2573
$ws.Range("I34").Value = 1206.5385
$ws.Range("J34").Value = 8603.817999999999
$ws.Range("K34").Value = 1206.5385
$ws.Range("L34").Value = 8603.817999999999
$ws.Range("M34").Value = -1004.5385
$ws.Range("N34").Value = -9007.817999999999
$ws.Range("H107").Value = 1109.3939
$ws.Range("I107").Value = 336.89474
$ws.Range("K107").Value = 336.89474
$ws.Range("M107").Value = 1583.10526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 143763.36
$ws.Range("I2").Value = 209.44444
$ws.Range("K2").Value = 1256.66664
$ws.Range("M2").Value = -1143.66664
$ws.Range("H107").Value = 1183.4
$ws.Range("J107").Value = 1926.7142
$ws.Range("L107").Value = 5780.142599999999
$ws.Range("N107").Value = -9620.142599999999
$ws.Range("H113").Value = 1668.9375
$ws.Range("I113").Value = 920.55554
$ws.Range("K113").Value = 2761.66662
$ws.Range("M113").Value = -591.66662
$ws.Range("H121").Value = 22224732
$ws.Range("J121").Value = 12502824
$ws.Range("L121").Value = 37508472
$ws.Range("N121").Value = -37511092
$ws.Range("H129").Value = 9863678
$ws.Range("I129").Value = 488.9
$ws.Range("J129").Value = 23953948
$ws.Range("K129").Value = 1466.7
$ws.Range("L129").Value = 71861844
$ws.Range("M129").Value = 3533.3
$ws.Range("N129").Value = -71871844
$ws.Range("H131").Value = 2039.9796
$ws.Range("I131").Value = 1718.5454
$ws.Range("J131").Value = 2133.0264
$ws.Range("K131").Value = 5155.6362
$ws.Range("L131").Value = 6399.0792
$ws.Range("M131").Value = -115.6361999999999
$ws.Range("N131").Value = -16479.0792
$ws.Range("H140").Value = 154626.16
$ws.Range("I140").Value = 167203.75
$ws.Range("J140").Value = 3695
$ws.Range("K140").Value = 501611.25
$ws.Range("L140").Value = 11085
$ws.Range("M140").Value = -496431.25
$ws.Range("N140").Value = -21445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H132").Value = 3581.7666
$ws.Range("I132").Value = 1923.6957
$ws.Range("J132").Value = 9029.714
$ws.Range("K132").Value = 5771.0871
$ws.Range("L132").Value = 27089.142
$ws.Range("M132").Value = -3241.0871
$ws.Range("N132").Value = -32149.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6170.32
$ws.Range("I7").Value = 5486.9287
$ws.Range("J7").Value = 7040.091
$ws.Range("K7").Value = 5486.9287
$ws.Range("L7").Value = 7040.091
$ws.Range("M7").Value = -5374.9287
$ws.Range("N7").Value = -7264.091
$ws.Range("H93").Value = 6202.476
$ws.Range("I93").Value = 6074.1177
$ws.Range("K93").Value = 6074.1177
$ws.Range("M93").Value = -4826.1177
$ws.Range("H126").Value = 6170.32
$ws.Range("I126").Value = 5486.9287
$ws.Range("J126").Value = 7040.091
$ws.Range("K126").Value = 16460.7861
$ws.Range("L126").Value = 21120.273
$ws.Range("M126").Value = -13990.7861
$ws.Range("N126").Value = -26060.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 7793.636
$ws.Range("J4").Value = 7533
$ws.Range("L4").Value = 7533
$ws.Range("N4").Value = -7759
$ws.Range("H53").Value = 25000000
$ws.Range("I53").Value = 25000000
$ws.Range("K53").Value = 25000000
$ws.Range("M53").Value = -24999393
